# VyTrackQa2Users.xlsx update
#
# The "storemanager52" test-user row (row 8) had its firstname/lastname
# placeholder replaced: "Eddie Rodriguez" -> "Roma Medhurst".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "Roma"
$ws.Range("D8").Value = "Medhurst"

# Column A (usernames) was widened so the longer values fit on screen.
$ws.Range("A:A").ColumnWidth = 26.14

# The active cell/selection moved to A8 before the file was saved.
[void]$ws.Range("A8").Select()
